$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws.Range("B3").Value = "6.0.0"

# Update Date value
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (was empty)
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it, shifting rows 12+ up
$ws.Rows.Item(11).Delete() | Out-Null

# After the delete, the "Case Sensitive" row (previously row 15) is now row 14; set its value to
# the literal text "true" (not the boolean TRUE) by computing it as a text formula, then
# converting the formula to a static value in place so the cell keeps its original style/type.
$c = $ws.Range("B14")
$c.Formula = '=TEXT("true","@")'
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false
